$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared-string entry holding the blank/placeholder "<Host ...>" XML
# snippet in B13 is being removed entirely (it was only ever blank text
# nodes wrapped in tags). After its removal, B13 takes on the value that
# was already used for the other "not configured" rows: "null".
$ws.Range("B13").Value = "null"
